$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 307.4737
$ws.Range("I33").Value = 242.46666
$ws.Range("J33").Value = 551.25
$ws.Range("K33").Value = 242.46666
$ws.Range("L33").Value = 551.25
$ws.Range("M33").Value = -13.46665999999999
$ws.Range("N33").Value = -1009.25

# Row 43
$ws.Range("H43").Value = 3545.4
$ws.Range("J43").Value = 3229.1428
$ws.Range("L43").Value = 3229.1428
$ws.Range("N43").Value = -3367.1428

# Row 86
$ws.Range("H86").Value = 7971
$ws.Range("I86").Value = 3000
$ws.Range("J86").Value = 8681.143
$ws.Range("K86").Value = 3000
$ws.Range("L86").Value = 8681.143
$ws.Range("M86").Value = -1877
$ws.Range("N86").Value = -10927.143

# Row 89
$ws.Range("H89").Value = 7971
$ws.Range("I89").Value = 3000
$ws.Range("J89").Value = 8681.143
$ws.Range("K89").Value = 15000
$ws.Range("L89").Value = 43405.715
$ws.Range("M89").Value = -9384
$ws.Range("N89").Value = -54637.715

# Row 138
$ws.Range("H138").Value = 3286.238
$ws.Range("I138").Value = 1700.8125
$ws.Range("K138").Value = 5102.4375
$ws.Range("M138").Value = 37.5625

# Row 141
$ws.Range("H141").Value = 16137119
$ws.Range("I141").Value = 22730888
$ws.Range("K141").Value = 68192664
$ws.Range("M141").Value = -68187484


$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2257.6667
$ws.Range("I32").Value = 2199.319
$ws.Range("K32").Value = 2199.319
$ws.Range("M32").Value = -1912.319

# Row 46
$ws.Range("H46").Value = 18232.727
$ws.Range("J46").Value = 19066.2
$ws.Range("L46").Value = 19066.2
$ws.Range("N46").Value = -19704.2

# Row 102
$ws.Range("H102").Value = 2736.2856
$ws.Range("I102").Value = 2139.2778
$ws.Range("J102").Value = 6318.3335
$ws.Range("K102").Value = 2139.2778
$ws.Range("L102").Value = 6318.3335
$ws.Range("M102").Value = -517.2777999999998
$ws.Range("N102").Value = -9562.333500000001

# Row 122
$ws.Range("H122").Value = 3377.7856
$ws.Range("I122").Value = 3253
$ws.Range("K122").Value = 9759
$ws.Range("M122").Value = -7309


$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 9270.182000000001
$ws.Range("I20").Value = 11386.625
$ws.Range("J20").Value = 3626.3333
$ws.Range("K20").Value = 11386.625
$ws.Range("L20").Value = 3626.3333
$ws.Range("M20").Value = -11139.625
$ws.Range("N20").Value = -4120.3333


$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 38857652
$ws.Range("I31").Value = 45458656
$ws.Range("J31").Value = 2552134.5
$ws.Range("K31").Value = 45458656
$ws.Range("L31").Value = 2552134.5
$ws.Range("M31").Value = -45458361
$ws.Range("N31").Value = -2552724.5

# Row 34
$ws.Range("H34").Value = 38857652
$ws.Range("I34").Value = 45458656
$ws.Range("J34").Value = 2552134.5
$ws.Range("K34").Value = 45458656
$ws.Range("L34").Value = 2552134.5
$ws.Range("M34").Value = -45458454
$ws.Range("N34").Value = -2552538.5

# Row 86
$ws.Range("H86").Value = 13158.6
$ws.Range("I86").Value = 13067.615
$ws.Range("K86").Value = 13067.615
$ws.Range("M86").Value = -11944.615

# Row 89
$ws.Range("H89").Value = 13158.6
$ws.Range("I89").Value = 13067.615
$ws.Range("K89").Value = 65338.075
$ws.Range("M89").Value = -59722.075

# Row 94
$ws.Range("H94").Value = 1537.25
$ws.Range("I94").Value = 1365.6666
$ws.Range("J94").Value = 1708.8334
$ws.Range("K94").Value = 1365.6666
$ws.Range("L94").Value = 1708.8334
$ws.Range("M94").Value = -914.6666
$ws.Range("N94").Value = -2610.8334

# Row 107
$ws.Range("H107").Value = 2698.4866
$ws.Range("I107").Value = 2494.074
$ws.Range("J107").Value = 3250.4
$ws.Range("K107").Value = 2494.074
$ws.Range("L107").Value = 3250.4
$ws.Range("M107").Value = -574.0740000000001
$ws.Range("N107").Value = -7090.4

# Row 132
$ws.Range("H132").Value = 2964.1875
$ws.Range("I132").Value = 2784.6667
$ws.Range("K132").Value = 8354.000100000001
$ws.Range("M132").Value = -5824.000100000001


$ws = $wb.Worksheets.Item("CUL")
# Row 11
$ws.Range("H11").Value = 7019.143
$ws.Range("I11").Value = 2560.2
$ws.Range("K11").Value = 7680.599999999999
$ws.Range("M11").Value = -7540.599999999999

# Row 26
$ws.Range("H26").Value = 612.25
$ws.Range("J26").Value = 500
$ws.Range("L26").Value = 1500
$ws.Range("N26").Value = -2076

# Row 86
$ws.Range("H86").Value = 270
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()

# Row 89
$ws.Range("H89").Value = 270
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()


$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 2993.2942
$ws.Range("I80").Value = 1802.8572
$ws.Range("J80").Value = 3826.6
$ws.Range("K80").Value = 1802.8572
$ws.Range("L80").Value = 3826.6
$ws.Range("M80").Value = -804.8571999999999
$ws.Range("N80").Value = -5822.6

# Row 83
$ws.Range("H83").Value = 2993.2942
$ws.Range("I83").Value = 1802.8572
$ws.Range("J83").Value = 3826.6
$ws.Range("K83").Value = 9014.286
$ws.Range("L83").Value = 19133
$ws.Range("M83").Value = -4022.286
$ws.Range("N83").Value = -29117

# Row 102
$ws.Range("H102").Value = 2145.28
$ws.Range("I102").Value = 2138.875
$ws.Range("K102").Value = 2138.875
$ws.Range("M102").Value = -516.875

# Row 113
$ws.Range("H113").Value = 928102.7
$ws.Range("I113").Value = 2136.7334
$ws.Range("J113").Value = 3706000.5
$ws.Range("K113").Value = 2136.7334
$ws.Range("L113").Value = 3706000.5
$ws.Range("M113").Value = 33.26659999999993
$ws.Range("N113").Value = -3710340.5


$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 4207.45
$ws.Range("I132").Value = 2989.1538
$ws.Range("K132").Value = 8967.4614
$ws.Range("M132").Value = -6437.4614


$ws = $wb.Worksheets.Item("WVR")
# Row 6
$ws.Range("H6").Value = 16481.334
$ws.Range("I6").Value = 44444
$ws.Range("J6").Value = 2500
$ws.Range("K6").Value = 44444
$ws.Range("L6").Value = 2500
$ws.Range("M6").Value = -44329
$ws.Range("N6").Value = -2730

# Row 41
$ws.Range("H41").Value = 29993.334
$ws.Range("I41").Value = 29993
$ws.Range("J41").Value = 29993.5
$ws.Range("K41").Value = 29993
$ws.Range("L41").Value = 29993.5
$ws.Range("N41").Value = -30773.5
$ws.Range("M41").Value = -29603

# Row 45
$ws.Range("H45").Value = 10702
$ws.Range("J45").Value = 6054.1113
$ws.Range("L45").Value = 6054.1113
$ws.Range("N45").Value = -7036.1113

